$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 22:
# value changes from 45208 to 45212 (serial date numbers)
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
